# Updated capital structure database
# Refresh the Hong Kong Retail (Automotive) comps: recompute metrics for the
# existing companies (rows 2-5), re-sort/relabel companies B3:B5, and drop
# Phoenitron Holdings Limited (row 6) which fell out of the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "3"
$ws.Range("D2").Value = -0.312
$ws.Range("G2").Value = -0.6802018260451707
$ws.Range("H2").Value = -0.6802018260451707
$ws.Range("I2").Value = -0.4887073522345027
$ws.Range("J2").Value = -0.4887073522345027
$ws.Range("K2").Value = -70.92
$ws.Range("L2").Value = -1.703988467083133
$ws.Range("U2").Value = 21.04
$ws.Range("V2").Value = 0.1379672131147541
$ws.Range("W2").Value = -116.3967611336032
$ws.Range("X2").Value = 0.05929950728461264
$ws.Range("Y2").Value = -116.4560606408878
$ws.Range("Z2").Value = 0.1091197013234893
$ws.Range("AA2").Value = -0.09634447497478348
$ws.Range("AB2").Value = 0.05424122330220779
$ws.Range("AC2").Value = -0.1505856982769913
$ws.Range("AD2").Value = 176.078
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 176.078
$ws.Range("AG2").Value = 155.038
$ws.Range("AH2").Value = 0.5358788476404386
$ws.Range("AI2").Value = 0.9214980269837448
$ws.Range("AJ2").Value = 0.5041263193491536
$ws.Range("AK2").Value = 0.9117844246580177
$ws.Range("AL2").Value = 15.794
$ws.Range("AM2").Value = 15.515
$ws.Range("AN2").Value = -9.512587790383577
$ws.Range("AO2").Value = -1.287830821831075
$ws.Range("AP2").Value = -8.375904916261481
$ws.Range("AQ2").Value = -1.310989365130519
$ws.Range("B3").Value = "AVIC Joy Holdings (HK) Limited (SEHK:260)"
$ws.Range("D3").Value = -0.605
$ws.Range("G3").Value = -4.348958333333333
$ws.Range("H3").Value = -4.348958333333333
$ws.Range("I3").Value = -3.859375
$ws.Range("J3").Value = -3.859375
$ws.Range("K3").Value = -57.5
$ws.Range("L3").Value = -29.94791666666667
$ws.Range("U3").Value = 14.8
$ws.Range("V3").Value = 0.6040816326530613
$ws.Range("W3").Value = -116.3967611336032
$ws.Range("X3").Value = 0.2824091471849559
$ws.Range("Y3").Value = -116.6791702807882
$ws.Range("Z3").Value = 0.005927499274498786
$ws.Range("AA3").Value = -0.02287644251251875
$ws.Range("AB3").Value = 0.061019645838779
$ws.Range("AC3").Value = -0.08389608835129775
$ws.Range("AD3").Value = 154.9
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 154.9
$ws.Range("AG3").Value = 140.1
$ws.Range("AH3").Value = 0.8634336677814939
$ws.Range("AI3").Value = 1.841854934601665
$ws.Range("AJ3").Value = 0.8511543134872418
$ws.Range("AK3").Value = 2.021645021645022
$ws.Range("AL3").Value = 14.6
$ws.Range("AM3").Value = 14.432
$ws.Range("AN3").Value = -21.97163120567376
$ws.Range("AO3").Value = -0.5075342465753425
$ws.Range("AP3").Value = -19.87234042553191
$ws.Range("AQ3").Value = -0.5134423503325942
$ws.Range("B4").Value = "Auto Italia Holdings Limited (SEHK:720)"
$ws.Range("D4").Value = -0.312
$ws.Range("G4").Value = -0.2594871794871795
$ws.Range("H4").Value = -0.2594871794871795
$ws.Range("I4").Value = -0.2841025641025641
$ws.Range("J4").Value = -0.2841025641025641
$ws.Range("K4").Value = -5.6
$ws.Range("L4").Value = -0.2871794871794872
$ws.Range("U4").Value = 6.24
$ws.Range("V4").Value = 0.05426086956521739
$ws.Range("W4").Value = -0.08945686900958466
$ws.Range("X4").Value = 0.05929950728461264
$ws.Range("Y4").Value = -0.1487563762941973
$ws.Range("Z4").Value = 0.3391186393516747
$ws.Range("AA4").Value = -0.09634447497478348
$ws.Range("AB4").Value = 0.05424122330220779
$ws.Range("AC4").Value = -0.1505856982769913
$ws.Range("AD4").Value = 20.6
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 20.6
$ws.Range("AG4").Value = 14.36
$ws.Range("AH4").Value = 0.1519174041297935
$ws.Range("AI4").Value = 0.2647814910025707
$ws.Range("AJ4").Value = 0.1110080395794681
$ws.Range("AK4").Value = 0.2006707657909447
$ws.Range("AL4").Value = 1.17
$ws.Range("AM4").Value = 1.109
$ws.Range("AN4").Value = -4.63963963963964
$ws.Range("AO4").Value = -4.735042735042735
$ws.Range("AP4").Value = -3.234234234234234
$ws.Range("AQ4").Value = -4.995491433724076
$ws.Range("B5").Value = "Sinofortune Financial Holdings Limited (SEHK:8123)"
$ws.Range("D5").Value = 0.339
$ws.Range("G5").Value = -0.7376237623762376
$ws.Range("H5").Value = -0.7376237623762376
$ws.Range("I5").Value = -0.3658415841584158
$ws.Range("J5").Value = -0.3658415841584158
$ws.Range("K5").Value = -7.82
$ws.Range("L5").Value = -0.3871287128712871
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").ClearContents()
$ws.Range("X5").Value = 0.0544086769060776
$ws.Range("Y5").ClearContents()
$ws.Range("Z5").ClearContents()
$ws.Range("AA5").ClearContents()
$ws.Range("AB5").Value = 0.05334806747920273
$ws.Range("AC5").ClearContents()
$ws.Range("AD5").Value = 0.578
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0.578
$ws.Range("AG5").Value = 0.578
$ws.Range("AH5").Value = 0.04256886139343055
$ws.Range("AI5").Value = 0.01980944547261635
$ws.Range("AJ5").Value = 0.04256886139343055
$ws.Range("AK5").Value = 0.01980944547261635
$ws.Range("AL5").Value = 0.024
$ws.Range("AM5").Value = -0.026
$ws.Range("AN5").Value = -0.08233618233618234
$ws.Range("AO5").Value = -307.9166666666666
$ws.Range("AP5").Value = -0.08233618233618234
$ws.Range("AQ5").Value = 284.2307692307692

# Remove the Phoenitron Holdings row entirely (data refresh dropped this company)
$ws.Rows("6:6").Delete()

